$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.400.61'
$ws.Range('E2').Value = '  +1.13%  '

$ws.Range('D3').Value = '2.913.19'
$ws.Range('E3').Value = '  +3.91%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('E5').Value = '  -0.43%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '112.67'
$ws.Range('E6').Value = '  +0.81%  '

$ws.Range('E7').Value = '  +0.80%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('E9').Value = '  -0.57%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.98'
$ws.Range('E10').Value = '  -0.98%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0866'
$ws.Range('E11').Value = '  +3.15%  '

$ws.Range('E12').Value = '  +0.21%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.89'
$ws.Range('E13').Value = '  -0.67%  '

$ws.Range('E14').Value = '  +0.63%  '

$ws.Range('D15').Value = '3.372.49'
$ws.Range('E15').Value = '  +3.96%  '

$ws.Range('E16').Value = '  +6.10%  '

$ws.Range('D17').Value = '2.908.55'
$ws.Range('E17').Value = '  +3.74%  '

$ws.Range('D18').Value = '52.416.28'
$ws.Range('E18').Value = '  +1.12%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.65'
$ws.Range('E19').Value = '  -0.26%  '

$ws.Range('E20').Value = '  +3.85%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.25'
$ws.Range('E21').Value = '  +3.84%  '

$ws.Range('D22').Value = '0.0₃0982'
$ws.Range('E22').Value = '  +0.30%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.92'
$ws.Range('E23').Value = '  +0.53%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '270.27'
$ws.Range('E24').Value = '  +0.58%  '

$ws.Range('E25').Value = '  +0.38%  '

$ws.Range('E26').Value = '  +8.08%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.80'
$ws.Range('E27').Value = '  +2.27%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.02%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.69'
$ws.Range('E29').Value = '  +2.68%  '

$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.21'
$ws.Range('E30').Value = '  -2.43%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.68'
$ws.Range('E31').Value = '  +8.59%  '

$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.101'
$ws.Range('E32').Value = '  +14.19%  '

$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.31'
$ws.Range('E33').Value = '  +11.87%  '

$ws.Range('E34').Value = '  +0.37%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '53.25'
$ws.Range('E35').Value = '  +1.79%  '

$ws.Range('E36').Value = '  +1.68%  '

$ws.Range('E37').Value = '  -0.07%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.36'
$ws.Range('E38').Value = '  +6.66%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.03'
$ws.Range('E39').Value = '  +0.66%  '

$ws.Range('E40').Value = '  +3.53%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.83'
$ws.Range('E41').Value = '  +12.72%  '

$ws.Range('E42').Value = '  +1.62%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.38'
$ws.Range('E43').Value = '  +6.17%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '121.13'
$ws.Range('E44').Value = '  +0.62%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.59'
$ws.Range('E45').Value = '  +6.93%  '

$ws.Range('E46').Value = '  -1.71%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.55'
$ws.Range('E47').Value = '  +3.63%  '

$ws.Range('D48').Value = '2.198.79'
$ws.Range('E48').Value = '  +3.79%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.261'
$ws.Range('E49').Value = '  +20.75%  '

$ws.Range('E50').Value = '  +12.08%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.965'
$ws.Range('E51').Value = '  +0.19%  '
